# reclamation-des-documents-de-fin-de-contrat.docx — typo fixes
#
# 1) "Lettre Recommandée avec Accusé de Réception"
#       -> "Lettre recommandée avec accusé de réception"
#    (only the leading letter of each word after the first stays
#     capitalised; the rest is lower-cased)
#
# 2) "l'expression" (straight apostrophe) -> "l’expression" (curly quote)

$d = $word.ActiveDocument

# -- 1. Fix capitalisation typo in the "Lettre recommandée ..." line ------
$d.Content.Find.Execute(
    "Lettre Recommandée avec Accusé de Réception",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Lettre recommandée avec accusé de réception", 2
)

# -- 2. Replace the straight apostrophe with a typographic one ------------
$d.Content.Find.Execute(
    "l'expression",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "l" + [char]0x2019 + "expression", 2
)
